$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Accounting/Currency number format matching existing style index 9 (numFmtId 44)
$currencyFormat = '_("$"* #,##0.00_);_("$"* \(#,##0.00\);_("$"* "-"??_);_(@_)'

# --- Recalculated ManagementCost results (utility mob-cost multiplier 10% -> 5%) ---
$ws.Range("G136").Value = 1883.71900612063
$ws.Range("H136").Value = 188371.90061206301
$ws.Range("I136").Value = 1.25581267074708
$ws.Range("G136:I136").NumberFormat = $currencyFormat
$ws.Range("G137").Value = 4007.2314022508899
$ws.Range("H137").Value = 400723.14022508898
$ws.Range("I137").Value = 2.6714876015005902
$ws.Range("G137:I137").NumberFormat = $currencyFormat
$ws.Range("G138").Value = 18877.634920373399
$ws.Range("H138").Value = 1887763.4920373401
$ws.Range("I138").Value = 12.5850899469156
$ws.Range("G138:I138").NumberFormat = $currencyFormat
$ws.Range("G139").Value = 3363.7839395011201
$ws.Range("H139").Value = 336378.39395011298
$ws.Range("I139").Value = 2.2425226263340798
$ws.Range("G139:I139").NumberFormat = $currencyFormat
$ws.Range("G140").Value = 43796.466892304699
$ws.Range("H140").Value = 4379646.68923047
$ws.Range("I140").Value = 29.1976445948698
$ws.Range("G140:I140").NumberFormat = $currencyFormat
$ws.Range("G141").Value = 20713.25
$ws.Range("H141").Value = 2071325
$ws.Range("I141").Value = 13.8088333333333
$ws.Range("G141:I141").NumberFormat = $currencyFormat
$ws.Range("G142").Value = 13067.4835
$ws.Range("H142").Value = 1306748.3500000001
$ws.Range("I142").Value = 8.71165566666666
$ws.Range("G142:I142").NumberFormat = $currencyFormat
$ws.Range("G169").Value = 2677.35219330378
$ws.Range("H169").Value = 160641.131598226
$ws.Range("I169").Value = 1.0709408773215101
$ws.Range("G169:I169").NumberFormat = $currencyFormat
$ws.Range("G170").Value = 6224.4272177893599
$ws.Range("H170").Value = 373465.63306736201
$ws.Range("I170").Value = 2.48977088711574
$ws.Range("G170:I170").NumberFormat = $currencyFormat
$ws.Range("G171").Value = 27236.217376004999
$ws.Range("H171").Value = 1634173.0425603001
$ws.Range("I171").Value = 10.894486950401999
$ws.Range("G171:I171").NumberFormat = $currencyFormat
$ws.Range("G172").Value = 4780.9860594710399
$ws.Range("H172").Value = 286859.16356826201
$ws.Range("I172").Value = 1.9123944237884101
$ws.Range("G172:I172").NumberFormat = $currencyFormat
$ws.Range("G173").Value = 62248.438494312897
$ws.Range("H173").Value = 3734906.30965877
$ws.Range("I173").Value = 24.899375397725098
$ws.Range("G173:I173").NumberFormat = $currencyFormat
$ws.Range("G174").Value = 29169.75
$ws.Range("H174").Value = 1750185
$ws.Range("I174").Value = 11.667899999999999
$ws.Range("G174:I174").NumberFormat = $currencyFormat
$ws.Range("G175").Value = 18172.713366666601
$ws.Range("H175").Value = 1090362.8019999999
$ws.Range("I175").Value = 7.2690853466666603
$ws.Range("G175:I175").NumberFormat = $currencyFormat
$ws.Range("G202").Value = 4484.0019930777999
$ws.Range("H202").Value = 125552.055806178
$ws.Range("I202").Value = 1.2455561091882801
$ws.Range("G202:I202").NumberFormat = $currencyFormat
$ws.Range("G203").Value = 11887.7056449095
$ws.Range("H203").Value = 332855.75805746799
$ws.Range("I203").Value = 3.3021404569193198
$ws.Range("G203:I203").NumberFormat = $currencyFormat
$ws.Range("G204").Value = 55117.699983321399
$ws.Range("H204").Value = 1543295.59953299
$ws.Range("I204").Value = 15.310472217589201
$ws.Range("G204:I204").NumberFormat = $currencyFormat
$ws.Range("G205").Value = 8007.1464162103703
$ws.Range("H205").Value = 224200.09965389001
$ws.Range("I205").Value = 2.2242073378362099
$ws.Range("G205:I205").NumberFormat = $currencyFormat
$ws.Range("G206").Value = 104253.046339059
$ws.Range("H206").Value = 2919085.2974936501
$ws.Range("I206").Value = 28.959179538627499
$ws.Range("G206:I206").NumberFormat = $currencyFormat
$ws.Range("G207").Value = 52491.1785714285
$ws.Range("H207").Value = 1469753
$ws.Range("I207").Value = 14.5808829365079
$ws.Range("G207:I207").NumberFormat = $currencyFormat
$ws.Range("G208").Value = 30052.964285714199
$ws.Range("H208").Value = 841483
$ws.Range("I208").Value = 8.3480456349206307
$ws.Range("G208:I208").NumberFormat = $currencyFormat
$ws.Range("G235").Value = 5342.9406247643301
$ws.Range("H235").Value = 149602.33749340099
$ws.Range("I235").Value = 1.4841501735456399
$ws.Range("G235:I235").NumberFormat = $currencyFormat
$ws.Range("G236").Value = 14302.8444649783
$ws.Range("H236").Value = 400479.64501939301
$ws.Range("I236").Value = 3.9730123513828701
$ws.Range("G236:I236").NumberFormat = $currencyFormat
$ws.Range("G237").Value = 55117.699983321399
$ws.Range("H237").Value = 1543295.59953299
$ws.Range("I237").Value = 15.310472217589201
$ws.Range("G237:I237").NumberFormat = $currencyFormat
$ws.Range("G238").Value = 9540.9654013648797
$ws.Range("H238").Value = 267147.03123821598
$ws.Range("I238").Value = 2.6502681670457999
$ws.Range("G238:I238").NumberFormat = $currencyFormat
$ws.Range("G239").Value = 124223.36952577
$ws.Range("H239").Value = 3478254.3467215798
$ws.Range("I239").Value = 34.506491534936302
$ws.Range("G239:I239").NumberFormat = $currencyFormat
$ws.Range("G240").Value = 60048.321428571398
$ws.Range("H240").Value = 1681353
$ws.Range("I240").Value = 16.6800892857142
$ws.Range("G240:I240").NumberFormat = $currencyFormat
$ws.Range("G241").Value = 30052.964285714199
$ws.Range("H241").Value = 841483
$ws.Range("I241").Value = 8.3480456349206307
$ws.Range("G241:I241").NumberFormat = $currencyFormat

# --- Switch the AutoFilter from FoundationCost+Mobilization to ManagementCost ---
$ws.AutoFilter.Range.AutoFilter(6)
$ws.AutoFilter.Range.AutoFilter(4, "ManagementCost")

# --- Update the saved selection/view state ---
$ws.Range("K139").Select()
